# Update the "想去人数" (number of people wanting to go) figures in column F
# across the relevant worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1331
$ws1.Range("F6").Value = 1702
$ws1.Range("F7").Value = 6218
$ws1.Range("F16").Value = 6927
$ws1.Range("F17").Value = 124
$ws1.Range("F26").Value = 1582
$ws1.Range("F27").Value = 758
$ws1.Range("F31").Value = 55
$ws1.Range("F33").Value = 3894

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2259

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2259
$ws4.Range("F10").Value = 1331
$ws4.Range("F12").Value = 1702
$ws4.Range("F13").Value = 6218
$ws4.Range("F23").Value = 6927
$ws4.Range("F24").Value = 124
$ws4.Range("F33").Value = 1582
$ws4.Range("F34").Value = 758
$ws4.Range("F45").Value = 3894
